# Update crypto price/volume figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "39.977.47"
Set-TextValue $ws.Range("E2") "  +2.53%  "
Set-TextValue $ws.Range("D3") "2.229.66"
Set-TextValue $ws.Range("E3") "  +0.85%  "
Set-TextValue $ws.Range("E4") "  +0.15%  "
Set-TextValue $ws.Range("D5") "293.67"
Set-TextValue $ws.Range("E5") "  -0.77%  "
Set-TextValue $ws.Range("D6") "85.65"
Set-TextValue $ws.Range("E6") "  +6.90%  "
Set-TextValue $ws.Range("D7") "0.513"
Set-TextValue $ws.Range("E7") "  +1.73%  "
Set-TextValue $ws.Range("E8") "  +0.03%  "
Set-TextValue $ws.Range("D9") "0.470"
Set-TextValue $ws.Range("E9") "  +2.94%  "
Set-TextValue $ws.Range("D10") "30.90"
Set-TextValue $ws.Range("E10") "  +10.55%  "
Set-TextValue $ws.Range("D11") "0.0786"
Set-TextValue $ws.Range("D12") "46.77"
Set-TextValue $ws.Range("E12") "  +1.75%  "
Set-TextValue $ws.Range("D13") "0.108"
Set-TextValue $ws.Range("E13") "  +1.08%  "
Set-TextValue $ws.Range("D14") "6.43"
Set-TextValue $ws.Range("E14") "  +5.58%  "
Set-TextValue $ws.Range("D15") "2.582.48"
Set-TextValue $ws.Range("E15") "  +1.04%  "
Set-TextValue $ws.Range("D16") "14.05"
Set-TextValue $ws.Range("E16") "  +1.19%  "
Set-TextValue $ws.Range("D17") "2.296.74"
Set-TextValue $ws.Range("E17") "  +3.31%  "
Set-TextValue $ws.Range("D18") "0.725"
Set-TextValue $ws.Range("E18") "  +2.24%  "
Set-TextValue $ws.Range("D19") "39.924.23"
Set-TextValue $ws.Range("E19") "  +2.63%  "
Set-TextValue $ws.Range("D20") "0.0₃0887"
Set-TextValue $ws.Range("E20") "  +3.49%  "
Set-TextValue $ws.Range("D21") "5.77"
Set-TextValue $ws.Range("E21") "  +1.71%  "
Set-TextValue $ws.Range("D22") "10.74"
Set-TextValue $ws.Range("E22") "  +9.81%  "
Set-TextValue $ws.Range("D23") "65.17"
Set-TextValue $ws.Range("E23") "  +0.93%  "
Set-TextValue $ws.Range("D24") "234.62"
Set-TextValue $ws.Range("E24") "  +4.56%  "
Set-TextValue $ws.Range("E25") "  -0.22%  "
Set-TextValue $ws.Range("D26") "2.45"
Set-TextValue $ws.Range("E26") "  +3.06%  "
Set-TextValue $ws.Range("E27") "  +5.18%  "
Set-TextValue $ws.Range("D28") "22.71"
Set-TextValue $ws.Range("E28") "  +2.34%  "
Set-TextValue $ws.Range("D29") "2.22"
Set-TextValue $ws.Range("E29") "  +2.68%  "
Set-TextValue $ws.Range("D30") "9.19"
Set-TextValue $ws.Range("E30") "  +3.41%  "
Set-TextValue $ws.Range("D31") "33.06"
Set-TextValue $ws.Range("E31") "  +6.51%  "
Set-TextValue $ws.Range("D32") "152.18"
Set-TextValue $ws.Range("E32") "  +2.75%  "
Set-TextValue $ws.Range("E33") "  +0.20%  "
Set-TextValue $ws.Range("D34") "4.85"
Set-TextValue $ws.Range("E34") "  +2.60%  "
Set-TextValue $ws.Range("D35") "0.0715"
Set-TextValue $ws.Range("E35") "  +4.34%  "
Set-TextValue $ws.Range("E36") "  +2.50%  "
Set-TextValue $ws.Range("D37") "16.17"
Set-TextValue $ws.Range("E37") "  +13.20%  "
Set-TextValue $ws.Range("E38") "  +2.53%  "
Set-TextValue $ws.Range("D39") "0.0994"
Set-TextValue $ws.Range("E39") "  +3.71%  "
Set-TextValue $ws.Range("D40") "2.70"
Set-TextValue $ws.Range("E40") "  +2.80%  "
Set-TextValue $ws.Range("D41") "1.68"
Set-TextValue $ws.Range("E41") "  +5.83%  "
Set-TextValue $ws.Range("D42") "3.80"
Set-TextValue $ws.Range("E42") "  +5.73%  "
Set-TextValue $ws.Range("D43") "2.034.23"
Set-TextValue $ws.Range("E43") "  +7.29%  "
Set-TextValue $ws.Range("D44") "2.21"
Set-TextValue $ws.Range("E44") "  +6.50%  "
Set-TextValue $ws.Range("D45") "0.0268"
Set-TextValue $ws.Range("E45") "  +5.85%  "
Set-TextValue $ws.Range("D46") "9.95"
Set-TextValue $ws.Range("E46") "  +13.31%  "
Set-TextValue $ws.Range("D47") "16.19"
Set-TextValue $ws.Range("E47") "  +0.55%  "
Set-TextValue $ws.Range("D48") "2.56"
Set-TextValue $ws.Range("E48") "  +2.26%  "
Set-TextValue $ws.Range("D49") "2.450.82"
Set-TextValue $ws.Range("E49") "  +1.26%  "
Set-TextValue $ws.Range("D50") "70.47"
Set-TextValue $ws.Range("E50") "  +1.32%  "
Set-TextValue $ws.Range("D51") "1.45"
Set-TextValue $ws.Range("E51") "  +14.02%  "
